# Correções de imports e conexão das interfaces
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Col A (name), Col B (code), Col C (label), Col D (number),
# Col E (bool/number or blank), Col F (bool or blank)
$data = @(
    @("auto1", "ar",    "A/C",       27,  $false, $null),
    @("auto2", "ar",    "A/C",       30,  $true,  $null),
    @("auto2", "tv",    "Televisor", 1,   100,    $true),
    @("auto2", "lamp1", "Lâmpada",   100, $true,  $null),
    @("AUTO3", "ar",    "A/C",       16,  $false, $null),
    @("AUTO3", "tv",    "Televisor", 1,   0,      $false),
    @("AUTO3", "lamp1", "Lâmpada",   0,   $false, $null),
    @("auto4", "ar",    "A/C",       30,  $true,  $null),
    @("auto4", "tv",    "Televisor", 5,   50,     $true),
    @("auto4", "lamp1", "Lâmpada",   52,  $true,  $null)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    if ($null -ne $row[4]) {
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
    if ($null -ne $row[5]) {
        $ws.Cells.Item($r, 6).Value = $row[5]
    }
}
